$d = $word.ActiveDocument

# Make sure edits apply directly (no tracked-changes markup)
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $false

# --- 1) Rewrite the "Reproducción, estrategia K, estrategia r" paragraph ---
# Locate it by its current text instead of a hard-coded index, so the script
# is resilient to any paragraph shift.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Reproducci*n, estrategia K, estrategia r*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the keywords paragraph"
}

$r = $target.Range

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES_tradnl"/></w:rPr>'

$frag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr>$rPr</w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r>$rPr<w:t>reproducci&#243;n,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>$rPr<w:t>estrategia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r>$rPr<w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r>$rPr<w:t>K,estrategia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r>$rPr<w:t xml:space="preserve"> r</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$r.InsertXML($frag)

# --- 2) Move the "_GoBack" bookmark from the end of the document to just
#        after the keywords paragraph ---
$afterIdx = $target.Index + 1
$bmRange = $d.Paragraphs.Item($afterIdx).Range
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.TrackRevisions = $wasTracking
